$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price record inserted at row 115, pushing the existing
# rows 115-223 down to 116-224 (dimension grows from T223 to T224).
$ws.Rows.Item(115).Insert()

$ws.Cells.Item(115, 1).Value = 4
$ws.Cells.Item(115, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(115, 3).Value = "Los Lagos"
$ws.Cells.Item(115, 4).Value = 44778
$ws.Cells.Item(115, 5).Value = 10
$ws.Cells.Item(115, 6).Value = "Fruta"
$ws.Cells.Item(115, 7).Value = 100108
$ws.Cells.Item(115, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(115, 9).Value = 100108002
$ws.Cells.Item(115, 10).Value = "Mango"
$ws.Cells.Item(115, 11).Value = "Sin especificar"
$ws.Cells.Item(115, 12).Value = "Primera"
$ws.Cells.Item(115, 13).Value = 120
$ws.Cells.Item(115, 14).Value = 13000
$ws.Cells.Item(115, 15).Value = 14000
$ws.Cells.Item(115, 16).Value = 13500
$ws.Cells.Item(115, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(115, 18).Value = "Brasil"
$ws.Cells.Item(115, 19).Value = 3375
$ws.Cells.Item(115, 20).Value = 4

# Keep the same date-number formatting as the other rows in column D.
$ws.Cells.Item(115, 4).NumberFormat = $ws.Cells.Item(116, 4).NumberFormat
